$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 30
$ws.Range("F3").Value = 1203
$ws.Range("F4").Value = 310
$ws.Range("F6").Value = 82
$ws.Range("F7").Value = 946
$ws.Range("F8").Value = 351
$ws.Range("F9").Value = 602
$ws.Range("F10").Value = 548
$ws.Range("F11").Value = 1419
$ws.Range("F12").Value = 131
$ws.Range("F13").Value = 1321
$ws.Range("F14").Value = 2976
$ws.Range("F15").Value = 380
$ws.Range("F16").Value = 1594
$ws.Range("F18").Value = 778
$ws.Range("F19").Value = 229
$ws.Range("F20").Value = 1350
$ws.Range("F21").Value = 256
$ws.Range("F22").Value = 62
$ws.Range("F24").Value = 392
$ws.Range("F25").Value = 3423
$ws.Range("F27").Value = 557
$ws.Range("F28").Value = 1521

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 46
$ws.Range("F8").Value = 19
$ws.Range("F12").Value = 71
$ws.Range("F13").Value = 14

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 791

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 30
$ws.Range("F3").Value = 791
$ws.Range("F4").Value = 1203
$ws.Range("F5").Value = 310
$ws.Range("F12").Value = 6
$ws.Range("F13").Value = 46
$ws.Range("F14").Value = 82
$ws.Range("F15").Value = 19
$ws.Range("F17").Value = 946
$ws.Range("F18").Value = 351
$ws.Range("F19").Value = 602
$ws.Range("F20").Value = 548
$ws.Range("F21").Value = 1419
$ws.Range("F22").Value = 131
$ws.Range("F23").Value = 1321
$ws.Range("F24").Value = 2976
$ws.Range("F25").Value = 380
$ws.Range("F26").Value = 1594
$ws.Range("F28").Value = 778
$ws.Range("F29").Value = 229
$ws.Range("F30").Value = 1350
$ws.Range("F31").Value = 256
$ws.Range("F32").Value = 62
$ws.Range("F36").Value = 392
$ws.Range("F37").Value = 3423
$ws.Range("F39").Value = 557
$ws.Range("F40").Value = 1521
$ws.Range("F41").Value = 71
$ws.Range("F42").Value = 14
